# Refresh the cryptocurrency price / 1h-volume figures (columns D and E)
# on the "cryptos" worksheet, per the scheduled GitHub Actions data pull.
#
# All values in columns D (Price) and E (Volume(1h)) are stored as TEXT in
# this sheet (e.g. "8.25", "  +16.97%  " with padding spaces) rather than
# numbers, so that values like "43.676.31" (thousand-separated) or strings
# with a literal trailing zero ("8.20") round-trip exactly.  To avoid Excel
# silently re-interpreting a numeric-looking string as a real Number (which
# would drop a significant trailing zero, e.g. "8.20" -> 8.2) we force the
# cell to Text format before writing the value, then restore the cell style
# to "Normal" afterwards so no stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, in the same order as the worksheet rows.
$updates = [ordered]@{
    'D2' = '43.541.76'
    'E2' = '  -1.38%  '
    'D3' = '2.226.48'
    'E3' = '  +0.29%  '
    'E4' = '  -0.03%  '
    'D5' = '269.78'
    'E5' = '  +3.75%  '
    'D6' = '92.44'
    'E6' = '  +11.59%  '
    'D7' = '0.623'
    'E7' = '  -1.11%  '
    'D9' = '0.619'
    'E9' = '  +2.02%  '
    'D10' = '45.96'
    'E10' = '  +4.30%  '
    'E11' = '  -0.74%  '
    'D12' = '8.20'
    'E12' = '  +16.04%  '
    'E13' = '  +0.88%  '
    'D14' = '2.563.64'
    'E14' = '  +0.35%  '
    'D15' = '15.05'
    'E15' = '  +3.06%  '
    'D16' = '2.227.67'
    'E16' = '  +1.07%  '
    'D17' = '0.801'
    'E17' = '  +3.02%  '
    'D18' = '43.513.54'
    'E18' = '  -1.24%  '
    'D19' = '0.0000103'
    'E19' = '  -0.69%  '
    'D20' = '5.99'
    'E20' = '  -0.38%  '
    'D21' = '70.34'
    'E21' = '  -1.37%  '
    'E22' = '  -1.61%  '
    'D23' = '232.85'
    'E23' = '  -0.10%  '
    'E24' = '  -3.86%  '
    'E25' = '  -0.04%  '
    'D26' = '2.51'
    'D27' = '11.31'
    'E27' = '  +4.64%  '
    'E28' = '  +5.42%  '
    'D29' = '40.57'
    'E29' = '  -2.30%  '
    'D30' = '2.26'
    'E30' = '  +1.95%  '
    'D31' = '172.44'
    'E31' = '  -0.28%  '
    'D32' = '0.0922'
    'E32' = '  +4.92%  '
    'D33' = '20.81'
    'E33' = '  +0.68%  '
    'D34' = '5.45'
    'E34' = '  +2.10%  '
    'E35' = '  +0.21%  '
    'E36' = '  -3.61%  '
    'E37' = '  -3.10%  '
    'D38' = '4.30'
    'E38' = '  -4.47%  '
    'D39' = '3.55'
    'E39' = '  +20.35%  '
    'D40' = '12.49'
    'E40' = '  -8.01%  '
    'D41' = '2.17'
    'E41' = '  +2.15%  '
    'D42' = '0.218'
    'E42' = '  +8.05%  '
    'D43' = '63.29'
    'E43' = '  -0.25%  '
    'D44' = '5.32'
    'E44' = '  -4.61%  '
    'D45' = '0.0987'
    'E45' = '  -0.23%  '
    'D46' = '8.37'
    'E46' = '  +0.58%  '
    'D47' = '100.30'
    'E47' = '  -2.82%  '
    'D48' = '1.15'
    'E48' = '  +2.51%  '
    'E49' = '  +1.39%  '
    'D50' = '0.438'
    'E50' = '  -1.03%  '
    'D51' = '2.452.71'
    'E51' = '  +0.61%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
